# Allow datetime objects in info dict
#
# Adds a "Test Date" column (B) to the experiment-record sheet: a text
# header in B1 and a date value in B2, formatted as a short date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell (new shared string "Test Date")
$ws.Range("B1").Value = "Test Date"

# New date value in B2 (serial 45588 == 2024-10-23), displayed as m/d/yyyy
# (built-in numFmtId 14). NumberFormatLocal applies the built-in date format
# directly instead of registering a redundant custom number format entry.
$ws.Range("B2").Value = 45588
$ws.Range("B2").NumberFormatLocal = "mm-dd-yy"

# Widen column B to fit the new "Test Date" header/date values.
$ws.Columns("B").ColumnWidth = 9.666666666666666

# Move the active selection to B3, just below the new data.
[void]$ws.Range("B3").Select()
